$wb = $excel.ActiveWorkbook

# Column F on each sheet is "想去人数" (want-to-go count).
# Values below are the new values taken from the target diff, keyed by row number.

$sheet1Updates = @{
    3  = 3389
    9  = 7
    10 = 1848
    11 = 423
    12 = 171
    13 = 1631
    14 = 59
    15 = 1098
    19 = 8549
    20 = 203
    21 = 1139
    22 = 312
    25 = 1023
    26 = 567
    27 = 1180
    29 = 591
    30 = 29
    33 = 112
    35 = 637
    36 = 471
    37 = 373
    38 = 3540
    39 = 937
    41 = 765
    43 = 114
    45 = 666
    46 = 47
    48 = 24
    49 = 2416
}

$sheet2Updates = @{
    10 = 203
    15 = 3
    17 = 21
    18 = 33
    19 = 16
    24 = 6942
    30 = 10
    32 = 19
}

$sheet3Updates = @{
    4  = 2074
    5  = 1397
    8  = 2237
    10 = 1316
}

$sheet4Updates = @{
    3  = 2074
    6  = 2237
    8  = 1316
    10 = 7
    11 = 423
    12 = 171
    14 = 8549
    15 = 203
    16 = 1139
    17 = 312
    19 = 1023
    20 = 1180
    21 = 591
    22 = 29
    25 = 203
    26 = 3
    27 = 112
    28 = 637
    29 = 373
    31 = 21
    32 = 33
    33 = 3540
    34 = 937
    40 = 666
    43 = 47
    45 = 24
    47 = 2416
}

function Apply-Updates($sheetName, $updates) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}

Apply-Updates "展览" $sheet1Updates
Apply-Updates "演出" $sheet2Updates
Apply-Updates "本地生活" $sheet3Updates
Apply-Updates "全部类型" $sheet4Updates

$wb.Save()
